$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("B19").Value = 6802951
$ws.Range("F19").Value = "Sparta Prague"
$ws.Range("G19").Value = "Pardubice"
$ws.Range("H19").Value = 5
$ws.Range("J19").Value = "H"
$ws.Range("K19").Value = 1.166
$ws.Range("L19").Value = 6
$ws.Range("M19").Value = 15
$ws.Range("N19").Value = 1.111
$ws.Range("O19").Value = 7.5
$ws.Range("P19").Value = 17
$ws.Range("Q19").Value = -2.25
$ws.Range("R19").Value = 1.85
$ws.Range("S19").Value = 2
$ws.Range("T19").Value = 3.5
$ws.Range("U19").Value = 1.825
$ws.Range("V19").Value = 2.025
$ws.Range("W19").Value = 0.111
$ws.Range("Y19").Value = -1
$ws.Range("Z19").Value = 0.8500000000000001
$ws.Range("AA19").Value = -1
$ws.Range("AB19").Value = 0.825
$ws.Range("AC19").Value = -1

# Row 20
$ws.Range("B20").Value = 6802956
$ws.Range("F20").Value = "Slovacko"
$ws.Range("G20").Value = "Sigma Olomouc"
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = "A"
$ws.Range("K20").Value = 1.909
$ws.Range("L20").Value = 3.4
$ws.Range("M20").Value = 3.75
$ws.Range("N20").Value = 1.75
$ws.Range("O20").Value = 3.5
$ws.Range("P20").Value = 4.5
$ws.Range("Q20").Value = -0.5
$ws.Range("R20").Value = 1.8
$ws.Range("S20").Value = 2.05
$ws.Range("T20").Value = 2.25
$ws.Range("U20").Value = 1.85
$ws.Range("V20").Value = 2
$ws.Range("W20").Value = -1
$ws.Range("Y20").Value = 3.5
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = 1.05
$ws.Range("AB20").Value = -0.5
$ws.Range("AC20").Value = 0.5

# Row 81
$ws.Range("B81").Value = 6803017
$ws.Range("F81").Value = "Bohemians 1905"
$ws.Range("G81").Value = "FC Trinity Zlin"
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = "D"
$ws.Range("K81").Value = 1.6
$ws.Range("L81").Value = 4
$ws.Range("M81").Value = 5.25
$ws.Range("N81").Value = 1.4
$ws.Range("O81").Value = 4.75
$ws.Range("P81").Value = 6.5
$ws.Range("Q81").Value = -1.25
$ws.Range("R81").Value = 1.9
$ws.Range("S81").Value = 1.95
$ws.Range("W81").Value = -1
$ws.Range("X81").Value = 3.75
$ws.Range("Z81").Value = -1
$ws.Range("AA81").Value = 0.95

# Row 82
$ws.Range("B82").Value = 6803015
$ws.Range("F82").Value = "Ceske Budejovice"
$ws.Range("G82").Value = "MFK Karvina"
$ws.Range("H82").Value = 1
$ws.Range("J82").Value = "H"
$ws.Range("K82").Value = 1.95
$ws.Range("L82").Value = 3.5
$ws.Range("M82").Value = 3.75
$ws.Range("N82").Value = 1.8
$ws.Range("O82").Value = 3.75
$ws.Range("P82").Value = 4
$ws.Range("Q82").Value = -0.5
$ws.Range("R82").Value = 1.85
$ws.Range("S82").Value = 2
$ws.Range("W82").Value = 0.8
$ws.Range("X82").Value = -1
$ws.Range("Z82").Value = 0.8500000000000001
$ws.Range("AA82").Value = -1

# Row 186
$ws.Range("B186").Value = 6803125
$ws.Range("F186").Value = "Slovacko"
$ws.Range("G186").Value = "FK Jablonec"
$ws.Range("H186").Value = 0
$ws.Range("I186").Value = 1
$ws.Range("J186").Value = "A"
$ws.Range("K186").Value = 1.8
$ws.Range("L186").Value = 3.4
$ws.Range("M186").Value = 4.2
$ws.Range("N186").Value = 1.7
$ws.Range("O186").Value = 3.5
$ws.Range("P186").Value = 4.5
$ws.Range("Q186").Value = -0.75
$ws.Range("R186").Value = 1.975
$ws.Range("S186").Value = 1.875
$ws.Range("T186").Value = 2.5
$ws.Range("U186").Value = 1.9
$ws.Range("V186").Value = 1.95
$ws.Range("W186").Value = -1
$ws.Range("Y186").Value = 3.5
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = 0.875
$ws.Range("AC186").Value = 0.95

# Row 187
$ws.Range("B187").Value = 6803127
$ws.Range("F187").Value = "Slovan Liberec"
$ws.Range("G187").Value = "MFK Karvina"
$ws.Range("H187").Value = 1
$ws.Range("K187").Value = 1.6
$ws.Range("L187").Value = 4
$ws.Range("M187").Value = 5
$ws.Range("N187").Value = 1.45
$ws.Range("O187").Value = 4.333
$ws.Range("P187").Value = 6
$ws.Range("Q187").Value = -1
$ws.Range("R187").Value = 1.875
$ws.Range("S187").Value = 1.975
$ws.Range("T187").Value = 2.75
$ws.Range("U187").Value = 1.85
$ws.Range("V187").Value = 2
$ws.Range("W187").Value = 0.45
$ws.Range("Z187").Value = 0
$ws.Range("AA187").Value = -0
$ws.Range("AC187").Value = 1

# Row 188
$ws.Range("B188").Value = 6803124
$ws.Range("F188").Value = "Hradec Kralove"
$ws.Range("G188").Value = "Pardubice"
$ws.Range("H188").Value = 2
$ws.Range("J188").Value = "H"
$ws.Range("K188").Value = 1.85
$ws.Range("M188").Value = 3.8
$ws.Range("N188").Value = 1.85
$ws.Range("O188").Value = 3.5
$ws.Range("P188").Value = 3.8
$ws.Range("Q188").Value = -0.5
$ws.Range("R188").Value = 2.05
$ws.Range("S188").Value = 1.8
$ws.Range("U188").Value = 1.9
$ws.Range("V188").Value = 1.95
$ws.Range("W188").Value = 0.8500000000000001
$ws.Range("X188").Value = -1
$ws.Range("Z188").Value = 1.05
$ws.Range("AA188").Value = -1
$ws.Range("AC188").Value = 0.95

# Row 189
$ws.Range("B189").Value = 6803122
$ws.Range("F189").Value = "Bohemians 1905"
$ws.Range("G189").Value = "Mlada Boleslav"
$ws.Range("I189").Value = 0
$ws.Range("J189").Value = "D"
$ws.Range("K189").Value = 2.25
$ws.Range("M189").Value = 2.875
$ws.Range("N189").Value = 2.2
$ws.Range("O189").Value = 3.4
$ws.Range("P189").Value = 3
$ws.Range("Q189").Value = -0.25
$ws.Range("U189").Value = 1.85
$ws.Range("V189").Value = 2
$ws.Range("X189").Value = 2.4
$ws.Range("Y189").Value = -1
$ws.Range("Z189").Value = -0.5
$ws.Range("AA189").Value = 0.4375
$ws.Range("AC189").Value = 1

# Row 214
$ws.Range("B214").Value = 6803145
$ws.Range("F214").Value = "Pardubice"
$ws.Range("G214").Value = "Ceske Budejovice"
$ws.Range("H214").Value = 1
$ws.Range("J214").Value = "D"
$ws.Range("K214").Value = 1.95
$ws.Range("L214").Value = 3.5
$ws.Range("M214").Value = 3.75
$ws.Range("N214").Value = 1.95
$ws.Range("O214").Value = 3.5
$ws.Range("P214").Value = 3.8
$ws.Range("Q214").Value = -0.5
$ws.Range("R214").Value = 1.925
$ws.Range("S214").Value = 1.925
$ws.Range("T214").Value = 2.75
$ws.Range("X214").Value = 2.5
$ws.Range("Y214").Value = -1
$ws.Range("AA214").Value = 0.925

# Row 216
$ws.Range("B216").Value = 6803144
$ws.Range("F216").Value = "FK Teplice"
$ws.Range("G216").Value = "Hradec Kralove"
$ws.Range("H216").Value = 0
$ws.Range("J216").Value = "A"
$ws.Range("K216").Value = 2.25
$ws.Range("L216").Value = 3.2
$ws.Range("M216").Value = 3.3
$ws.Range("N216").Value = 2.1
$ws.Range("O216").Value = 3.3
$ws.Range("P216").Value = 3.6
$ws.Range("Q216").Value = -0.25
$ws.Range("R216").Value = 1.8
$ws.Range("S216").Value = 2.05
$ws.Range("T216").Value = 2.5
$ws.Range("X216").Value = -1
$ws.Range("Y216").Value = 2.6
$ws.Range("AA216").Value = 1.05

# Row 218
$ws.Range("B218").Value = 6803151
$ws.Range("F218").Value = "Bohemians 1905"
$ws.Range("G218").Value = "Sparta Prague"
$ws.Range("K218").Value = 5.25
$ws.Range("L218").Value = 4.333
$ws.Range("M218").Value = 1.571
$ws.Range("N218").Value = 5.25
$ws.Range("O218").Value = 4.333
$ws.Range("P218").Value = 1.571
$ws.Range("Q218").Value = 1
$ws.Range("R218").Value = 1.875
$ws.Range("S218").Value = 1.975
$ws.Range("T218").Value = 2.75
$ws.Range("U218").Value = 1.925
$ws.Range("V218").Value = 1.925

# Row 219
$ws.Range("R219").Value = 1.925
$ws.Range("S219").Value = 1.925

# Row 220
$ws.Range("B220").Value = 6851033
$ws.Range("F220").Value = "Mlada Boleslav"
$ws.Range("G220").Value = "FK Teplice"
$ws.Range("K220").Value = 1.75
$ws.Range("L220").Value = 3.75
$ws.Range("M220").Value = 4.5
$ws.Range("N220").Value = 1.666
$ws.Range("O220").Value = 3.8
$ws.Range("P220").Value = 5
$ws.Range("Q220").Value = -0.75
$ws.Range("R220").Value = 1.825
$ws.Range("S220").Value = 2.025
$ws.Range("U220").Value = 2
$ws.Range("V220").Value = 1.85

# Row 223
$ws.Range("N223").Value = 2.2
$ws.Range("O223").Value = 3.25
$ws.Range("P223").Value = 3.3
$ws.Range("R223").Value = 1.925
$ws.Range("S223").Value = 1.925

# Row 225
$ws.Range("R225").Value = 1.875
$ws.Range("S225").Value = 1.975
